# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described in the commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A known plain-text styled cell (header-free "Normal" style) used to strip
# any NumberFormat overrides we temporarily apply below, so cell styling stays untouched.
$normalStyle = $ws.Range("B9").Style

$ws.Range("D2").Value = '37.790.92'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.047.40'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.04%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '227.88'
$c.Style = $normalStyle
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("E6").Value = '  -0.50%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '59.84'
$c.Style = $normalStyle
$ws.Range("E7").Value = '  +0.46%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.377'
$c.Style = $normalStyle
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  +2.90%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = '2.351.88'
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("E13").Value = '  -1.66%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '21.43'
$c.Style = $normalStyle
$ws.Range("E14").Value = '  +1.36%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '5.50'
$c.Style = $normalStyle
$ws.Range("E15").Value = '  +6.23%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.764'
$c.Style = $normalStyle
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = '2.045.54'
$ws.Range("E17").Value = '  +0.83%  '
$ws.Range("D18").Value = '37.783.62'
$ws.Range("E18").Value = '  +0.05%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '69.52'
$c.Style = $normalStyle
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("E20").Value = '  -1.71%  '
$ws.Range("D21").Value = '0.0₃0829'
$ws.Range("E21").Value = '  +0.54%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '222.19'
$c.Style = $normalStyle
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("E25").Value = '  +2.96%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '168.99'
$c.Style = $normalStyle
$ws.Range("E26").Value = '  +2.45%  '
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  -0.97%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '18.78'
$c.Style = $normalStyle
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("E31").Value = '  -0.48%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.Style = $normalStyle
$ws.Range("E32").Value = '  +8.13%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '4.38'
$c.Style = $normalStyle
$ws.Range("E33").Value = '  -1.10%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.54'
$c.Style = $normalStyle
$ws.Range("E34").Value = '  +0.89%  '
$ws.Range("E35").Value = '  +0.06%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '6.48'
$c.Style = $normalStyle
$ws.Range("E36").Value = '  +1.91%  '
$ws.Range("E37").Value = '  +4.18%  '
$ws.Range("E38").Value = '  +7.84%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = $normalStyle
$ws.Range("E39").Value = '  -0.18%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.43'
$c.Style = $normalStyle
$ws.Range("E40").Value = '  +9.69%  '
$ws.Range("D41").Value = '1.524.45'
$ws.Range("E41").Value = '  +0.01%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '97.59'
$c.Style = $normalStyle
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("E43").Value = '  -1.25%  '
$ws.Range("E44").Value = '  -0.13%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '4.21'
$c.Style = $normalStyle
$ws.Range("E45").Value = '  +2.49%  '
$ws.Range("E46").Value = '  -2.94%  '
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '2.94'
$c.Style = $normalStyle
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '7.10'
$c.Style = $normalStyle
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").Value = '2.240.11'
$ws.Range("E51").Value = '  +0.92%  '
